# Insert a new "is_targeted list" worksheet (TRUE/FALSE) right before the
# "resolution_x_unit list" sheet, and repoint column N's data validation at
# it instead of the old inline "TRUE,FALSE" literal list.

$wb = $excel.ActiveWorkbook

# --- 1. Insert the new sheet in the right tab position ------------------
$refSheet = $wb.Worksheets.Item("resolution_x_unit list")
$newSheet = $wb.Worksheets.Add($refSheet)
$newSheet.Name = "is_targeted list"

# --- 2. Populate it with literal text "TRUE" / "FALSE" ------------------
# Plain Range.Value assignment of the strings "TRUE"/"FALSE" gets silently
# coerced to real booleans by Excel, which is not what the source list
# sheet wants (it needs literal text so the dropdown / validation formula
# compares against text). Route the literal text through a formula +
# paste-special-values round trip on a scratch cell, which keeps it a
# plain string cell.
$scratch = $wb.Worksheets.Item(1)
$scratchCell = $scratch.Range("ZZ1")

$scratchCell.Formula = '="TRUE"'
$scratchCell.Copy()
$newSheet.Range("A1").PasteSpecial(-4163)

$scratchCell.Formula = '="FALSE"'
$scratchCell.Copy()
$newSheet.Range("A2").PasteSpecial(-4163)

$scratchCell.ClearContents()
$excel.CutCopyMode = $false

# --- 3. Point column N's validation at the new list sheet ----------------
$ws = $wb.Worksheets.Item("Export as TSV")
$dv = $ws.Range("N2:N1048576").Validation
$dv.Formula1 = "='is_targeted list'!`$A`$1:`$A`$2"
$dv.ErrorTitle = "Value must come from list"
$dv.ErrorMessage = "Value must be one of: TRUE / FALSE."
